$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("données12")

# --- Data corrections (logic-problem fixes) ---
# Row 12
$ws.Range("A12").Value = 39.72
$ws.Range("C12").Value = 111

# Row 27
$ws.Range("A27").Value = 22.93
$ws.Range("C27").Value = 115

# Row 33
$ws.Range("A33").Value = 24.709999999999997
$ws.Range("C33").Value = 113

# --- Cosmetic / metadata changes (best effort) ---
# Resize the workbook window to match the new saved view geometry.
$excel.ActiveWindow.WindowState = -4137
$excel.ActiveWindow.Width = 25800
$excel.ActiveWindow.Height = 13200
